# Daily attendance processing - swap "Recorded By" value order
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in column G (Recorded By) wherever this exact combined value occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$updated = 0
for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
        $updated = $updated + 1
    }
}

Write-Output "Updated cells: $updated"
